# Generate Report for Handback
# Update the timestamp strings recorded on the Overview, zh-cn and de-de
# sheets to reflect the latest handoff/handback generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the
#     1e13c8ca-...-edf2cd4be70e.md row (row 3) ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-15 18:41:47"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H3) and
#     "Correspond Handback DateTime" (K3) for the 1e13c8ca-... row ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-15 18:41:42"
$wsZhCn.Range("K3").Value = "2016-08-15 18:41:59"

# --- de-de sheet: "Correspond Handback DateTime" (K3) for the
#     1e13c8ca-... row ---
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-15 18:42:13"
